# "finish group by job" -- turn the vertical list of job names (A1:A10)
# into a horizontal header row (A1:J1) and turn the whole former
# single-column list into an empty 15-col x 10-row grid ready for
# grouping party members by job underneath each header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 10 job names, in the same order they used to appear top-to-bottom
# in column A; they now become the column headers across row 1.
$jobs = @("奶", "火", "圣骑", "拳", "弩", "船", "饺子", "刀", "弓", "单挂")

# Build the new 15 (A:O) column x 10 row grid. Every cell keeps the big
# 18pt header-style font (style index 1 in the original file); row 1
# gets the job names across columns A-J as headers, everything else is
# blank (but still styled) space for the grouped members.
for ($r = 1; $r -le 10; $r++) {
    for ($c = 1; $c -le 15; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Font.Size = 18
        if ($r -eq 1 -and $c -le 10) {
            $cell.Value = $jobs[$c - 1]
        } else {
            $cell.Value = ""
        }
    }
}

# Uniform column widths across the new grid (~20.6 characters wide).
$ws.Range("A1:O10").ColumnWidth = 19.666666666666668

# Update the remembered selection to roughly where editing finished.
$ws.Range("L13").Select()
